$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets(1)
$ws.Range("F3").Value = 131
$ws.Range("F5").Value = 3
$ws.Range("F7").Value = 13046
$ws.Range("F8").Value = 57
$ws.Range("F10").Value = 261
$ws.Range("F11").Value = 3087
$ws.Range("F12").Value = 82
$ws.Range("F13").Value = 6499
$ws.Range("F16").Value = 3434
$ws.Range("F17").Value = 165
$ws.Range("F18").Value = 123
$ws.Range("F24").Value = 3617
$ws.Range("F27").Value = 2799
$ws.Range("F29").Value = 1889
$ws.Range("F30").Value = 103
$ws.Range("F32").Value = 6653
$ws.Range("F34").Value = 1006
$ws.Range("F35").Value = 1987
$ws.Range("F38").Value = 1042
$ws.Range("F40").Value = 211
$ws.Range("F41").Value = 222
$ws.Range("F44").Value = 138
$ws.Range("F45").Value = 1206
$ws.Range("F46").Value = 1787

$ws = $wb.Worksheets(2)
$ws.Range("F2").Value = 47
$ws.Range("F14").Value = 99

$ws = $wb.Worksheets(3)
$ws.Range("F2").Value = 438
$ws.Range("F3").Value = 604

$ws = $wb.Worksheets(4)
$ws.Range("F3").Value = 131
$ws.Range("F6").Value = 438
$ws.Range("F7").Value = 604
$ws.Range("F9").Value = 13046
$ws.Range("F10").Value = 57
$ws.Range("F13").Value = 261
$ws.Range("F14").Value = 3087
$ws.Range("F16").Value = 3434
$ws.Range("F17").Value = 165
$ws.Range("F24").Value = 3617
$ws.Range("F27").Value = 2799
$ws.Range("F28").Value = 2799
$ws.Range("F30").Value = 1889
$ws.Range("F31").Value = 103
$ws.Range("F33").Value = 6653
$ws.Range("F34").Value = 99
$ws.Range("F36").Value = 1006
$ws.Range("F37").Value = 1987
$ws.Range("F41").Value = 1042
$ws.Range("F42").Value = 211
$ws.Range("F43").Value = 222
$ws.Range("F45").Value = 1206
$ws.Range("F47").Value = 1787
